# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps that get refreshed each time the
# handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 8898bf84-...
$wsOverview.Range("G2").Value = "2016-08-31 11:12:28"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 8898bf84-...
$wsZhCn.Range("H2").Value = "2016-08-31 11:12:23"
$wsZhCn.Range("K2").Value = "2016-08-31 11:12:40"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 8898bf84-...
$wsDeDe.Range("H2").Value = "2016-08-31 11:12:28"
$wsDeDe.Range("K2").Value = "2016-08-31 11:12:47"
